$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "E3"  = 16.466
    "B7"  = 5.153
    "A10" = -21.604
    "A12" = -21.606
    "B15" = 4.867
    "A18" = -21.914
    "E18" = 16.284
    "E19" = 16.483
    "B20" = 6.406999999999999
    "E27" = 16.328
    "B29" = 5.237
    "B30" = 6.02
    "B31" = 5.847
    "A37" = -19.92
    "B40" = 9.327999999999999
    "E42" = 16.586
    "E44" = 16.753
    "E47" = 16.278
    "A55" = -21.795
    "E58" = 16.567
    "A68" = -21.473
    "B68" = 5.512
    "E73" = 16.572
    "B76" = 5.81
    "A77" = -20.899
    "A78" = -20.309
    "B87" = 4.836
    "B88" = 5.058000000000001
    "E95" = 17.564
    "B96" = 6.692
    "B98" = 5.095000000000001
    "B101" = 7.628000000000002
    "E101" = 16.677
    "B102" = 7.74
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
